$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("J2").Value = 8
$ws.Range("K2").Value = 4
$ws.Range("AE2").Value = $True
$ws.Range("AF2").Value = "05/27/2024 22:27:28"
$ws.Range("AG2").Value = "05/27/2024 23:26:44"
$ws.Range("AH2").Value = 0.3933170545426222

# Row 3
$ws.Range("J3").Value = 8
$ws.Range("K3").Value = 4
$ws.Range("AE3").Value = $True
$ws.Range("AF3").Value = "05/27/2024 23:26:44"
$ws.Range("AG3").Value = "05/28/2024 01:30:01"
$ws.Range("AH3").Value = 0.3891914388824514

# Row 4
$ws.Range("J4").Value = 8
$ws.Range("K4").Value = 4
$ws.Range("AE4").Value = $True
$ws.Range("AF4").Value = "05/28/2024 01:30:01"
$ws.Range("AG4").Value = "05/28/2024 03:01:17"
$ws.Range("AH4").Value = 0.382039836203902

# Row 5
$ws.Range("J5").Value = 8
$ws.Range("K5").Value = 4
$ws.Range("AE5").Value = $True
$ws.Range("AF5").Value = "05/28/2024 03:01:17"
$ws.Range("AG5").Value = "05/28/2024 04:07:27"
$ws.Range("AH5").Value = 0.4360767771293618

# Row 6
$ws.Range("J6").Value = 8
$ws.Range("K6").Value = 4
$ws.Range("AE6").Value = $True
$ws.Range("AF6").Value = "05/28/2024 04:07:27"
$ws.Range("AG6").Value = "05/28/2024 05:29:52"
$ws.Range("AH6").Value = 0.4189910766361019

# Row 7
$ws.Range("J7").Value = 8
$ws.Range("K7").Value = 4
$ws.Range("AE7").Value = $True
$ws.Range("AF7").Value = "05/28/2024 05:29:52"
$ws.Range("AG7").Value = "05/28/2024 19:10:59"
$ws.Range("AH7").Value = 0.3209438879973027

# Row 8
$ws.Range("J8").Value = 8
$ws.Range("K8").Value = 4
$ws.Range("AE8").Value = $True
$ws.Range("AF8").Value = "05/28/2024 05:29:53"
$ws.Range("AG8").Value = "05/28/2024 19:10:60"
$ws.Range("AH8").Value = 1.3209438879973

# Row 9
$ws.Range("J9").Value = 8
$ws.Range("K9").Value = 4
$ws.Range("AE9").Value = $True
$ws.Range("AF9").Value = "05/28/2024 05:29:54"
$ws.Range("AG9").Value = "05/28/2024 19:10:61"
$ws.Range("AH9").Value = 2.3209438879973

# Row 10
$ws.Range("J10").Value = 8
$ws.Range("K10").Value = 4
$ws.Range("AE10").Value = $True
$ws.Range("AF10").Value = "05/28/2024 05:29:55"
$ws.Range("AG10").Value = "05/28/2024 19:10:62"
$ws.Range("AH10").Value = 3.3209438879973

# Row 11
$ws.Range("J11").Value = 8
$ws.Range("K11").Value = 4
$ws.Range("AE11").Value = $True
$ws.Range("AF11").Value = "05/28/2024 05:29:56"
$ws.Range("AG11").Value = "05/28/2024 19:10:63"
$ws.Range("AH11").Value = 4.3209438879973

# Row 12
$ws.Range("J12").Value = 8
$ws.Range("K12").Value = 4
$ws.Range("AE12").Value = $True
$ws.Range("AF12").Value = "05/28/2024 19:28:25"
$ws.Range("AG12").Value = "05/28/2024 23:05:36"
$ws.Range("AH12").Value = 0.3244014611680931

# Row 13
$ws.Range("J13").Value = 8
$ws.Range("K13").Value = 4
$ws.Range("AE13").Value = $True
$ws.Range("AF13").Value = "05/28/2024 23:05:36"
$ws.Range("AG13").Value = "05/29/2024 02:56:58"
$ws.Range("AH13").Value = 0.3297904134435905

# Row 14
$ws.Range("J14").Value = 8
$ws.Range("K14").Value = 4
$ws.Range("AE14").Value = $True
$ws.Range("AF14").Value = "05/29/2024 03:12:29"
$ws.Range("AG14").Value = "05/29/2024 06:44:01"
$ws.Range("AH14").Value = 0.3689008988453189

# Row 15
$ws.Range("J15").Value = 8
$ws.Range("K15").Value = 4

# Row 16
$ws.Range("J16").Value = 8
$ws.Range("K16").Value = 4

# Row 17
$ws.Range("J17").Value = 2
$ws.Range("K17").Value = 1
$ws.Range("AE17").Value = $True
$ws.Range("AF17").Value = "05/28/2024 23:05:36"
$ws.Range("AG17").Value = "05/29/2024 02:56:58"
$ws.Range("AH17").Value = 0.3297904134435905

# Row 18
$ws.Range("J18").Value = 2
$ws.Range("K18").Value = 1
$ws.Range("AE18").Value = $True
$ws.Range("AF18").Value = "05/28/2024 23:05:36"
$ws.Range("AG18").Value = "05/29/2024 02:56:58"
$ws.Range("AH18").Value = 0.3297904134435905

# Row 19
$ws.Range("J19").Value = 2
$ws.Range("K19").Value = 1
$ws.Range("AE19").Value = $True
$ws.Range("AF19").Value = "05/28/2024 23:05:36"
$ws.Range("AG19").Value = "05/29/2024 02:56:58"
$ws.Range("AH19").Value = 0.3297904134435905

# Row 20
$ws.Range("J20").Value = 2
$ws.Range("K20").Value = 1
$ws.Range("AE20").Value = $True
$ws.Range("AF20").Value = "05/28/2024 23:05:36"
$ws.Range("AG20").Value = "05/29/2024 02:56:58"
$ws.Range("AH20").Value = 0.3297904134435905

# Row 21
$ws.Range("J21").Value = 2
$ws.Range("K21").Value = 1
$ws.Range("AE21").Value = $True
$ws.Range("AF21").Value = "05/28/2024 23:05:36"
$ws.Range("AG21").Value = "05/29/2024 02:56:58"
$ws.Range("AH21").Value = 0.3297904134435905

